# Apply the 2025-2 regular timetable update ('25.07.20 data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Targeted cell corrections (schedule / room / remark updates) that do NOT
#    involve any row shifting. All of these rows are well above the inserted
#    row (1414), so they are unaffected by the later row insert.
# ---------------------------------------------------------------------------

# 지역사회간호학II (간호학과) - team-taught sections now have a schedule/room
$ws.Range("J78").Value = "화5~6,금5~6"
$ws.Range("K78").Value = "다니엘관404호(大)강의실"

$ws.Range("J79").Value = "화7~8,금7~8"
$ws.Range("K79").Value = "다니엘관405호(大)강의실"

# 채플 (약학과 6년제)
$ws.Range("J693").Value = "월3"
$ws.Range("K693").Value = "미지정"

# 의약품분석학Ⅱ (약학과 6년제) - now scheduled
$ws.Range("J713").Value = "월6~7"
$ws.Range("K713").Value = "다니엘관303호(中)강의실"

# 약용식물학* (약학과 6년제) - schedule removed
$ws.Range("J714").Value = ""
$ws.Range("K714").Value = ""

# 국제경영학 (경영학과)
$ws.Range("J778").Value = "화1~3"
$ws.Range("K778").Value = "제3과학관418호(약학강의실)"

# 커뮤니케이션 전략 (경영학과)
$ws.Range("J784").Value = "월6~8"
$ws.Range("K784").Value = "제3과학관418호(약학강의실)"

# 교육방법 및 교육공학 (공통교양, 교직필수)
$ws.Range("J818").Value = "월5~7"
$ws.Range("K818").Value = "바울관312호강의실"
$ws.Range("M818").Value = "교직필수"

# 교육봉사활동 (공통교양, 교직필수) - schedule removed, remark added
$ws.Range("J819").Value = ""
$ws.Range("K819").Value = ""
$ws.Range("M819").Value = "교직필수"

# 교육심리 (공통교양, 교직필수) - schedule removed, remark added
$ws.Range("J820").Value = ""
$ws.Range("K820").Value = ""
$ws.Range("M820").Value = "교직필수"

# (KLEC)대학한국어 읽기Ⅱ - professor placeholder removed
$ws.Range("I842").Value = ""

# AI를 활용한 미래 도시문화 콘텐츠 - remark added
$ws.Range("M845").Value = "사회과학영역"

# 녹색도시와 인간생활 - remark added
$ws.Range("M917").Value = "사회과학영역"

# 대학한국어 말하기Ⅱ - professor removed, schedule/room added
$ws.Range("I923").Value = ""
$ws.Range("J923").Value = "수2~4"
$ws.Range("K923").Value = "바울관214호강의실"

# 대학한국어 쓰기Ⅱ - professor removed, schedule/room added
$ws.Range("I924").Value = ""
$ws.Range("J924").Value = "목6~9"
$ws.Range("K924").Value = "미지정"

# 문학적 상상력의 3D 구현
$ws.Range("J934").Value = "월5~7"
$ws.Range("K934").Value = "미지정"

# 문화와 사회, 그리고 인간(문사인(文社人)) - schedule/room added
$ws.Range("J935").Value = "목5~7"
$ws.Range("K935").Value = "미지정"

# 요한복음으로 살피는 예수의 메시지
$ws.Range("J991").Value = "수5~7"
$ws.Range("K991").Value = "미지정"

# 인공지능과 건축
$ws.Range("J996").Value = "수6~8"
$ws.Range("K996").Value = "미지정"

# 인문-디지털 가상현실의 문자적 경험
$ws.Range("J997").Value = "월2~4"
$ws.Range("K997").Value = "미지정"

# 재미와 상식으로 살펴보는 인류와 질병의 뒷이야기 - remark added
$ws.Range("M1001").Value = "자연과학영역"

# 치유적 연극: 나를 찾는 무대 - professor removed, remark added
$ws.Range("I1015").Value = ""
$ws.Range("M1015").Value = "인문예술영역"

# 탄소중립과 미래기술
$ws.Range("J1020").Value = "화6~8"
$ws.Range("K1020").Value = "미지정"

# ---------------------------------------------------------------------------
# 2) New course row inserted at row 1414 (항공관광외국어학부 - 일본문화와 언어),
#    pushing every following row (old 1414-1445) down by one (new 1415-1446).
# ---------------------------------------------------------------------------
$ws.Rows.Item(1414).Insert()

# The sheet stores every value (including numeric-looking codes) as text, so
# force the Text number format before writing - otherwise the COM layer's
# Value setter auto-coerces digit strings like "1815" into real numbers.
$ws.Range("A1414:N1414").NumberFormat = "@"

$ws.Range("A1414").Value = "1815"
$ws.Range("B1414").Value = "1004176"
$ws.Range("C1414").Value = "일본문화와 언어"
$ws.Range("D1414").Value = "항공관광외국어학부"
$ws.Range("E1414").Value = "1"
$ws.Range("F1414").Value = "전공선택"
$ws.Range("G1414").Value = ""
$ws.Range("H1414").Value = "3"
$ws.Range("I1414").Value = ""
$ws.Range("J1414").Value = ""
$ws.Range("K1414").Value = ""
$ws.Range("L1414").Value = "창의융합대학"
$ws.Range("M1414").Value = ""
$ws.Range("N1414").Value = ""
